$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 5767.3
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 5767.3
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 5767.3
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -8013.3
$ws.Range("H89").Value = 5767.3
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 5767.3
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 28836.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -40068.5
$ws.Range("H96").Value = 6031.722
$ws.Range("I96").Value = 8922.083000000001
$ws.Range("J96").Value = 251
$ws.Range("K96").Value = 26766.249
$ws.Range("L96").Value = 753
$ws.Range("M96").Value = -25393.249
$ws.Range("N96").Value = -3499
$ws.Range("H98").Value = 3398.1875
$ws.Range("I98").Value = 1187.5
$ws.Range("J98").Value = 3714
$ws.Range("K98").Value = 1187.5
$ws.Range("L98").Value = 3714
$ws.Range("M98").Value = 310.5
$ws.Range("N98").Value = -6710
$ws.Range("H122").Value = 3398.1875
$ws.Range("I122").Value = 1187.5
$ws.Range("J122").Value = 3714
$ws.Range("K122").Value = 3562.5
$ws.Range("L122").Value = 11142
$ws.Range("M122").Value = -1112.5
$ws.Range("N122").Value = -16042
$ws.Range("H125").Value = 375000350
$ws.Range("I125").Value = 500000160
$ws.Range("J125").Value = 250000500
$ws.Range("K125").Value = 4500001440
$ws.Range("L125").Value = 2250004500
$ws.Range("M125").Value = -4499998980
$ws.Range("N125").Value = -2250009420
$ws.Range("H132").Value = 1091.6
$ws.Range("I132").Value = 1067
$ws.Range("J132").Value = 1190
$ws.Range("K132").Value = 3201
$ws.Range("L132").Value = 3570
$ws.Range("M132").Value = -671
$ws.Range("N132").Value = -8630

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5363.9287
$ws.Range("I32").Value = 5007.577
$ws.Range("J32").Value = 9996.5
$ws.Range("K32").Value = 5007.577
$ws.Range("L32").Value = 9996.5
$ws.Range("M32").Value = -4720.577
$ws.Range("N32").Value = -10570.5
$ws.Range("H33").Value = 10013
$ws.Range("I33").Value = 10013
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 10013
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -9684
$ws.Range("H61").Value = 4401.857
$ws.Range("I61").Value = 2008.5
$ws.Range("J61").Value = 5874.6924
$ws.Range("K61").Value = 2008.5
$ws.Range("L61").Value = 5874.6924
$ws.Range("M61").Value = -1796.5
$ws.Range("N61").Value = -6298.6924
$ws.Range("H136").Value = 4401.857
$ws.Range("I136").Value = 2008.5
$ws.Range("J136").Value = 5874.6924
$ws.Range("K136").Value = 6025.5
$ws.Range("L136").Value = 17624.0772
$ws.Range("M136").Value = -3475.5
$ws.Range("N136").Value = -22724.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 210.75
$ws.Range("I64").Value = 97.5
$ws.Range("J64").Value = 248.5
$ws.Range("K64").Value = 97.5
$ws.Range("L64").Value = 248.5
$ws.Range("M64").Value = 127.5
$ws.Range("N64").Value = -698.5
$ws.Range("H67").Value = 210.75
$ws.Range("I67").Value = 97.5
$ws.Range("J67").Value = 248.5
$ws.Range("K67").Value = 97.5
$ws.Range("L67").Value = 248.5
$ws.Range("M67").Value = 682.5
$ws.Range("N67").Value = -1808.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2948.7144
$ws.Range("I58").Value = 3128.2
$ws.Range("J58").Value = 2500
$ws.Range("K58").Value = 3128.2
$ws.Range("L58").Value = 2500
$ws.Range("M58").Value = -2925.2
$ws.Range("N58").Value = -2906
$ws.Range("H136").Value = 2948.7144
$ws.Range("I136").Value = 3128.2
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 9384.599999999999
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -6834.599999999999
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100609.41
$ws.Range("I2").Value = 84709.92
$ws.Range("J2").Value = 123575.336
$ws.Range("K2").Value = 508259.52
$ws.Range("L2").Value = 741452.0159999999
$ws.Range("M2").Value = -508146.52
$ws.Range("N2").Value = -741678.0159999999
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 5
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 15
$ws.Range("M9").ClearContents()
$ws.Range("N9").Value = -463
$ws.Range("H12").Value = 30.11111
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 42.666668
$ws.Range("K12").Value = 15
$ws.Range("L12").Value = 128.000004
$ws.Range("M12").Value = 158
$ws.Range("N12").Value = -474.000004
$ws.Range("H92").Value = 228.47058
$ws.Range("I92").Value = 160.08333
$ws.Range("J92").Value = 392.6
$ws.Range("K92").Value = 480.24999
$ws.Range("L92").Value = 1177.8
$ws.Range("M92").Value = 767.75001
$ws.Range("N92").Value = -3673.8
$ws.Range("H93").Value = 20000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 20000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 60000
$ws.Range("M93").Value = -63744
$ws.Range("H107").Value = 101.5
$ws.Range("I107").Value = 103
$ws.Range("J107").Value = 100
$ws.Range("K107").Value = 309
$ws.Range("L107").Value = 300
$ws.Range("M107").Value = 1611
$ws.Range("N107").Value = -4140
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").ClearContents()
$ws.Range("H132").Value = 5299.125
$ws.Range("I132").Value = 5878.8
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 52909.2
$ws.Range("L132").Value = 38997
$ws.Range("M132").Value = -50379.2
$ws.Range("N132").Value = -44057
$ws.Range("H134").Value = 7184.8887
$ws.Range("I134").Value = 1627.3334
$ws.Range("J134").Value = 18300
$ws.Range("K134").Value = 4882.0002
$ws.Range("L134").Value = 54900
$ws.Range("M134").Value = 187.9997999999996
$ws.Range("N134").Value = -65040

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6926.385
$ws.Range("I7").Value = 2300.1428
$ws.Range("J7").Value = 8630.789000000001
$ws.Range("K7").Value = 2300.1428
$ws.Range("L7").Value = 8630.789000000001
$ws.Range("M7").Value = -2188.1428
$ws.Range("N7").Value = -8854.789000000001
$ws.Range("H22").Value = 1599.5
$ws.Range("I22").Value = 998
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 1599.5
$ws.Range("I27").Value = 998
$ws.Range("H46").Value = 35871.465
$ws.Range("I46").Value = 73439
$ws.Range("J46").Value = 2999.875
$ws.Range("K46").Value = 73439
$ws.Range("L46").Value = 2999.875
$ws.Range("M46").Value = -73251
$ws.Range("N46").Value = -3375.875
$ws.Range("H126").Value = 6926.385
$ws.Range("I126").Value = 2300.1428
$ws.Range("J126").Value = 8630.789000000001
$ws.Range("K126").Value = 6900.428400000001
$ws.Range("L126").Value = 25892.367
$ws.Range("M126").Value = -4430.428400000001
$ws.Range("N126").Value = -30832.367
$ws.Range("H136").Value = 2919.6667
$ws.Range("I136").Value = 2412.7778
$ws.Range("J136").Value = 3426.5557
$ws.Range("K136").Value = 7238.3334
$ws.Range("L136").Value = 10279.6671
$ws.Range("M136").Value = -4688.3334
$ws.Range("N136").Value = -15379.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1025000
$ws.Range("I2").Value = 2000000
$ws.Range("J2").Value = 700000
$ws.Range("K2").Value = 2000000
$ws.Range("L2").Value = 700000
$ws.Range("M2").Value = -1999888
$ws.Range("N2").Value = -700224
$ws.Range("H62").Value = 14945.5
$ws.Range("I62").Value = 14799.5
$ws.Range("J62").Value = 14982
$ws.Range("K62").Value = 14799.5
$ws.Range("L62").Value = 14982
$ws.Range("M62").Value = -14175.5
$ws.Range("N62").Value = -16230
$ws.Range("H65").Value = 14945.5
$ws.Range("I65").Value = 14799.5
$ws.Range("J65").Value = 14982
$ws.Range("K65").Value = 73997.5
$ws.Range("L65").Value = 74910
$ws.Range("M65").Value = -70877.5
$ws.Range("N65").Value = -81150
$ws.Range("H100").Value = 10002680
$ws.Range("I100").Value = 11112978
$ws.Range("J100").Value = 10000
$ws.Range("K100").Value = 22225956
$ws.Range("L100").Value = 20000
$ws.Range("M100").Value = -22225415
$ws.Range("N100").Value = -21082
$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 7500
$ws.Range("N122").Value = -12400
